# Regenerate save_data to use K (column G) instead of Strike# : recalculated
# K values (std/mean recomputed externally) are written back into column G
# for data rows 2-57 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (col G) for rows 2..57, in order.
$kValues = @(
    1,1,1,1,0,3,1,1,0,2,
    1,1,0,1,1,0,1,1,2,1,
    3,2,0,1,2,0,0,2,1,1,
    2,1,2,1,1,0,2,0,0,0,
    0,1,1,2,1,0,1,0,0,0,
    0,0,1,0,0,0
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
